$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the "Description_date" column (B): it was stored with an errant
#     leading minus sign on every record; replace each with its absolute
#     (positive) value. Row 1 is the text header, so start at row 2. ---
for ($r = 2; $r -le 86; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $v = $cell.Value2
    if ($v -ne $null) {
        $cell.Value = [Math]::Abs($v)
    }
}

# --- The old sheet had a stray empty-but-underline-formatted cell at N8
#     left over from editing; Clear() removes both its (nonexistent)
#     content and its formatting, which also shrinks the sheet's used
#     range/dimension back down from column N to column K. ---
[void]$ws.Range("N8").Clear()

# --- Re-point the saved view: select columns L:AA (now past the data,
#     matching the sheet as last saved) with L1 as the active cell, and
#     scroll the window so column E is the first visible column. ---
[void]$ws.Range("L1:AA1048576").Select()
$excel.ActiveWindow.ScrollColumn = 5
